$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Force text format on the changed price cells so values like '38.60' or
# '0.999' are preserved as text strings (matching the original inlineStr
# cells) rather than being auto-converted to floating point numbers.
$dCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D14', 'D15', 'D17', 'D18', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D32', 'D34', 'D39', 'D41', 'D42', 'D45', 'D46', 'D47', 'D49', 'D51')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '43.996.63'
$ws.Range('D3').Value = '2.264.56'
$ws.Range('D5').Value = '318.85'
$ws.Range('D6').Value = '102.69'
$ws.Range('D7').Value = '0.586'
$ws.Range('D9').Value = '0.570'
$ws.Range('D10').Value = '38.60'
$ws.Range('D12').Value = '7.86'
$ws.Range('D14').Value = '2.610.22'
$ws.Range('D15').Value = '0.874'
$ws.Range('D17').Value = '2.263.29'
$ws.Range('D18').Value = '43.908.90'
$ws.Range('D22').Value = '66.07'
$ws.Range('D23').Value = '3.20'
$ws.Range('D24').Value = '239.19'
$ws.Range('D25').Value = '2.20'
$ws.Range('D26').Value = '0.999'
$ws.Range('D28').Value = '10.25'
$ws.Range('D29').Value = '38.66'
$ws.Range('D32').Value = '163.22'
$ws.Range('D34').Value = '20.48'
$ws.Range('D39').Value = '4.54'
$ws.Range('D41').Value = '3.88'
$ws.Range('D42').Value = '15.80'
$ws.Range('D45').Value = '1.779.21'
$ws.Range('D46').Value = '0.208'
$ws.Range('D47').Value = '85.28'
$ws.Range('D49').Value = '8.92'
$ws.Range('D51').Value = '74.63'

# --- Volume(1h) (column E) updates ---
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E6').Value = '  +3.16%  '
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  +5.40%  '
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('E17').Value = '  +1.90%  '
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  +14.38%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('E31').Value = '  +2.10%  '
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('E41').Value = '  +7.78%  '
$ws.Range('E42').Value = '  +31.08%  '
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('E47').Value = '  -3.15%  '
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('E49').Value = '  +4.20%  '
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('E51').Value = '  -3.78%  '
